$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 9 to make room for the new entry "WEBサイトのスクレイピング",
# shifting the former rows 9-15 down to rows 10-16.
$ws.Rows.Item(9).Insert()

# Drop all existing hyperlink relationships; they will be rebuilt below so that every
# F-column cells hyperlink target matches its (possibly shifted) displayed URL exactly.
$ws.Hyperlinks.Delete()

$timestamp = "2025-09-11 18:20:39"

# Row 2
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【AIで開発生産性を革新】AI活用推進エンジニア募集(副業・業務委託)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5391864", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391864")
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# Row 3
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【業務委託/副業可】AI SaaS開発を牽引するCTO候補を募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5391872", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391872")
$ws.Range("G3").Value = 375
$ws.Range("H3").Value = "🔥AI,Ai ◆開発"

# Row 4
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【業務委託/副業可】AI SaaS開発を牽引するCTO候補を募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5391756", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391756")
$ws.Range("G4").Value = 375
$ws.Range("H4").Value = "🔥AI,Ai ◆開発"

# Row 5
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "【AIで開発生産性を革新】AI活用推進エンジニア募集(副業・業務委託)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5391761", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391761")
$ws.Range("G5").Value = 375
$ws.Range("H5").Value = "🔥AI,Ai ◆開発"

# Row 6
$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "競馬AIの開発ができる方、もしくはすでに開発済みの方"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5391744", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391744")
$ws.Range("G6").Value = 375
$ws.Range("H6").Value = "🔥AI,Ai ◆開発"

# Row 7
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "【AI技術顧問/戦略アドバイザー募集】最先端AIで事業の非連続な成長を牽引するエキスパート求む"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5391776", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391776")
$ws.Range("G7").Value = 310
$ws.Range("H7").Value = "🔥AI,Ai"

# Row 8
$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "【日本人限定/継続案件】Node.jsエンジニア募集(スクレイピング機能開発)"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5391607", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391607")
$ws.Range("G8").Value = 155
$ws.Range("H8").Value = "◆開発,Node.js"

# Row 9
$ws.Range("A9").Value = $timestamp
$ws.Range("B9").Value = "WEBサイトのスクレイピング"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5392043", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5392043")
$ws.Range("G9").Value = 65
$ws.Range("H9").Value = "◆スクレイピング ◇サイト"

# Row 10
$ws.Range("A10").Value = $timestamp
$ws.Range("B10").Value = "【急募】SharePoint+Power Platformでの不動産賃貸管理システム構築"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5391490", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391490")
$ws.Range("G10").Value = 60
$ws.Range("H10").Value = "◇管理"

# Row 11
$ws.Range("A11").Value = $timestamp
$ws.Range("B11").Value = "【急募】Salesforce・MA・CRMコンサルタント経験者を探しています!"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5371747", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5371747")
$ws.Range("G11").Value = 48
$ws.Range("H11").Value = "◆コンサル"

# Row 12
$ws.Range("A12").Value = $timestamp
$ws.Range("B12").Value = "【急募・再掲載】自社アプリのデバッグ・バグチェック業務依頼 ※NDA締結必須"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "~ 5,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5391844", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391844")
$ws.Range("G12").Value = 30
$ws.Range("H12").Value = "◇アプリ"

# Row 13
$ws.Range("A13").Value = $timestamp
$ws.Range("B13").Value = "初回 【フルリモート】フリーランスエンジニア募集"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5391489", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391489")
$ws.Range("G13").Value = 25

# Row 14
$ws.Range("A14").Value = $timestamp
$ws.Range("B14").Value = "要件定義や基本設計ができる方(1人月、約2年アサイン予定)"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5391221", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391221")
$ws.Range("G14").Value = 25

# Row 15
$ws.Range("A15").Value = $timestamp
$ws.Range("B15").Value = "【講師募集】Gensparkを使ったWEB構築チュートリアル募集"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5390165", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5390165")
$ws.Range("G15").Value = 13

# Row 16
$ws.Range("A16").Value = $timestamp
$ws.Range("B16").Value = "【急募】Googleアナリティクス連携の専門家を探しています"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5391267", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5391267")
$ws.Range("G16").Value = 10

# Column H was widened from 13 to 15 characters
$ws.Columns.Item(8).ColumnWidth = 15
